# Trade #98 closed at 2026-02-16 21:39:20 - leadlag UP +0.000%
#
# This script:
#  1. Marks the existing OPEN leadlag trade (row 51, Trade #62) as CLOSED
#     with its exit price / P&L / exit reason / duration filled in.
#  2. Appends a copy of that now-closed trade to the "All Trades" sheet.
#  3. Appends a brand-new OPEN trade (Trade #98) to the "leadlag" sheet.
#  4. Refreshes the aggregate stats on "Summary" and "Comparison".
#
# NOTE: values such as "66.1%", "2.83" or "2026-02-16" must be written as
# literal text (matching the source workbook's inlineStr cells), so a
# leading apostrophe is used to stop Excel from auto-converting them to
# numbers/dates/percentages.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) leadlag sheet: close out Trade #62 (row 51)
# ---------------------------------------------------------------
$leadlag = $wb.Worksheets.Item("leadlag")

$leadlag.Range("G51").Value = 68924.156124
$leadlag.Range("H51").Value = "CLOSED"
$leadlag.Range("I51").Value = -0.2976
$leadlag.Range("J51").Value = -2.98
$leadlag.Range("M51").Value = "time_exit_5min"
$leadlag.Range("N51").Value = 5

# ---------------------------------------------------------------
# 2) leadlag sheet: append new OPEN Trade #98 as row 74
# ---------------------------------------------------------------
$leadlag.Range("A74").Value = 98
$leadlag.Range("B74").Value = "'2026-02-16"
$leadlag.Range("C74").Value = "21:39:20"
$leadlag.Range("D74").Value = "leadlag"
$leadlag.Range("E74").Value = "UP"
$leadlag.Range("F74").Value = 68407.14999999999
$leadlag.Range("H74").Value = "OPEN"
$leadlag.Range("I74").Value = 0
$leadlag.Range("J74").Value = 0
$leadlag.Range("K74").Value = 0.75
$leadlag.Range("L74").Value = "Binance leading with 0.172% move"
$leadlag.Range("N74").Value = 0

# ---------------------------------------------------------------
# 3) All Trades sheet: append the now-closed trade as row 63
# ---------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("A63").Value = 62
$allTrades.Range("B63").Value = "'2026-02-16"
$allTrades.Range("C63").Value = "21:34:14"
$allTrades.Range("D63").Value = "leadlag"
$allTrades.Range("E63").Value = "DOWN"
$allTrades.Range("F63").Value = 68719.61500000001
$allTrades.Range("G63").Value = 68924.156124
$allTrades.Range("H63").Value = "CLOSED"
$allTrades.Range("I63").Value = -0.2976
$allTrades.Range("J63").Value = -2.98
$allTrades.Range("K63").Value = 0.7332
$allTrades.Range("L63").Value = "Coinbase leading with -0.073% move"
$allTrades.Range("M63").Value = "time_exit_5min"
$allTrades.Range("N63").Value = 5

# ---------------------------------------------------------------
# 4) Summary sheet: refresh OVERALL and leadlag aggregate rows
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("C2").Value = 62
$summary.Range("D2").Value = "'66.1%"
$summary.Range("E2").Value = "'+15.6799%"
$summary.Range("F2").Value = "'+0.2529%"

$summary.Range("C3").Value = 72
$summary.Range("D3").Value = "'43.1%"
$summary.Range("E3").Value = "'+10.8291%"
$summary.Range("F3").Value = "'+0.1504%"

# ---------------------------------------------------------------
# 5) Comparison sheet: refresh leadlag row
# ---------------------------------------------------------------
$comparison = $wb.Worksheets.Item("Comparison")

$comparison.Range("B2").Value = 72
$comparison.Range("C2").Value = "'43.1%"
$comparison.Range("D2").Value = "'2.83"
$comparison.Range("F2").Value = "'-0.3119%"
